$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 838
$ws.Cells.Item(7, 6).Value = 420
$ws.Cells.Item(8, 6).Value = 4700
$ws.Cells.Item(9, 6).Value = 4700
$ws.Cells.Item(12, 6).Value = 157
$ws.Cells.Item(16, 6).Value = 7492
$ws.Cells.Item(21, 6).Value = 522
$ws.Cells.Item(22, 6).Value = 1365
$ws.Cells.Item(25, 6).Value = 1747
$ws.Cells.Item(28, 6).Value = 6176
$ws.Cells.Item(29, 6).Value = 142
$ws.Cells.Item(30, 6).Value = 22
$ws.Cells.Item(31, 6).Value = 117
$ws.Cells.Item(34, 6).Value = 6420
$ws.Cells.Item(39, 6).Value = 20
$ws.Cells.Item(42, 6).Value = 30
$ws.Cells.Item(43, 6).Value = 60
$ws.Cells.Item(46, 6).Value = 439
$ws.Cells.Item(47, 6).Value = 2144
$ws.Cells.Item(49, 6).Value = 1077

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 36
$ws.Cells.Item(10, 6).Value = 8

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(6, 6).Value = 73
$ws.Cells.Item(8, 6).Value = 420
$ws.Cells.Item(9, 6).Value = 4700
$ws.Cells.Item(10, 6).Value = 4700
$ws.Cells.Item(13, 6).Value = 157
$ws.Cells.Item(17, 6).Value = 7492
$ws.Cells.Item(20, 6).Value = 522
$ws.Cells.Item(21, 6).Value = 1365
$ws.Cells.Item(24, 6).Value = 1747
$ws.Cells.Item(25, 6).Value = 36
$ws.Cells.Item(29, 6).Value = 6176
$ws.Cells.Item(30, 6).Value = 142
$ws.Cells.Item(31, 6).Value = 8
$ws.Cells.Item(32, 6).Value = 22
$ws.Cells.Item(33, 6).Value = 117
$ws.Cells.Item(36, 6).Value = 6420
$ws.Cells.Item(40, 6).Value = 20
$ws.Cells.Item(43, 6).Value = 30
$ws.Cells.Item(46, 6).Value = 439
$ws.Cells.Item(48, 6).Value = 2144
$ws.Cells.Item(49, 6).Value = 45
